$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1937716262975779
$ws.Range("C2").Value = 0.5813148788927336
$ws.Range("J2").Value = 0.02422145328719723
$ws.Range("P2").Value = 0.1418685121107267
$ws.Range("S2").Value = 0.05882352941176471
$ws.Range("B3").Value = 0.01136363636363636
$ws.Range("C3").Value = 0.02272727272727273
$ws.Range("J3").Value = 0.03409090909090909
$ws.Range("P3").Value = 0.7613636363636364
$ws.Range("S3").Value = 0.1704545454545454
$ws.Range("P4").Value = 0.8125
$ws.Range("S4").Value = 0.1875
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.05729166666666666
$ws.Range("D6").Value = 0.015625
$ws.Range("F6").Value = 0.05208333333333334
$ws.Range("J6").Value = 0.3020833333333333
$ws.Range("O6").Value = 0.046875
$ws.Range("Q6").Value = 0.1510416666666667
$ws.Range("R6").Value = 0.1041666666666667
$ws.Range("S6").Value = 0.2708333333333333
$ws.Range("B7").Value = 0.07731958762886598
$ws.Range("D7").Value = 0.01030927835051546
$ws.Range("F7").Value = 0.02577319587628866
$ws.Range("J7").Value = 0.154639175257732
$ws.Range("O7").Value = 0.0154639175257732
$ws.Range("Q7").Value = 0.2061855670103093
$ws.Range("R7").Value = 0.1030927835051546
$ws.Range("S7").Value = 0.4072164948453608
$ws.Range("B8").Value = 0.06635071090047394
$ws.Range("D8").Value = 0.01658767772511848
$ws.Range("E8").Value = 0.002369668246445498
$ws.Range("F8").Value = 0.05924170616113744
$ws.Range("J8").Value = 0.0995260663507109
$ws.Range("O8").Value = 0.01421800947867299
$ws.Range("Q8").Value = 0.2132701421800948
$ws.Range("R8").Value = 0.1137440758293839
$ws.Range("S8").Value = 0.4146919431279621
$ws.Range("B9").Value = 0.1075949367088608
$ws.Range("D9").Value = 0.02531645569620253
$ws.Range("E9").Value = 0.006329113924050633
$ws.Range("F9").Value = 0.04430379746835443
$ws.Range("J9").Value = 0.1329113924050633
$ws.Range("O9").Value = 0.0189873417721519
$ws.Range("Q9").Value = 0.1962025316455696
$ws.Range("R9").Value = 0.1075949367088608
$ws.Range("S9").Value = 0.3607594936708861
$ws.Range("B10").Value = 0.1170909090909091
$ws.Range("D10").Value = 0.01309090909090909
$ws.Range("F10").Value = 0.056
$ws.Range("J10").Value = 0.1272727272727273
$ws.Range("O10").Value = 0.01236363636363636
$ws.Range("Q10").Value = 0.2378181818181818
$ws.Range("R10").Value = 0.1098181818181818
$ws.Range("S10").Value = 0.3265454545454545
$ws.Range("G11").Value = 0.1451612903225807
$ws.Range("J11").Value = 0.1064516129032258
$ws.Range("K11").Value = 0.2032258064516129
$ws.Range("L11").Value = 0.5451612903225806
$ws.Range("G12").Value = 0.7705882352941177
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.005882352941176471
$ws.Range("L12").Value = 0.01764705882352941
$ws.Range("S12").Value = 0.005882352941176471
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("S13").Value = 0.05405405405405406
$ws.Range("F15").Value = 0.02136752136752137
$ws.Range("H15").Value = 0.1837606837606838
$ws.Range("I15").Value = 0.07264957264957266
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.05128205128205128
$ws.Range("M15").Value = 0.01282051282051282
$ws.Range("N15").Value = 0.004273504273504274
$ws.Range("O15").Value = 0.07264957264957266
$ws.Range("S15").Value = 0.2478632478632479
$ws.Range("F16").Value = 0.015625
$ws.Range("H16").Value = 0.2552083333333333
$ws.Range("I16").Value = 0.04166666666666666
$ws.Range("J16").Value = 0.3854166666666667
$ws.Range("K16").Value = 0.09895833333333333
$ws.Range("M16").Value = 0.02604166666666667
$ws.Range("N16").Value = 0.005208333333333333
$ws.Range("O16").Value = 0.08333333333333333
$ws.Range("S16").Value = 0.08854166666666667
$ws.Range("F17").Value = 0.01730769230769231
$ws.Range("H17").Value = 0.1480769230769231
$ws.Range("I17").Value = 0.06346153846153846
$ws.Range("J17").Value = 0.4769230769230769
$ws.Range("K17").Value = 0.1211538461538462
$ws.Range("M17").Value = 0.01730769230769231
$ws.Range("O17").Value = 0.05384615384615385
$ws.Range("S17").Value = 0.1019230769230769
$ws.Range("F18").Value = 0.02755905511811024
$ws.Range("H18").Value = 0.1732283464566929
$ws.Range("I18").Value = 0.06299212598425197
$ws.Range("J18").Value = 0.5236220472440944
$ws.Range("K18").Value = 0.07874015748031496
$ws.Range("M18").Value = 0.01574803149606299
$ws.Range("O18").Value = 0.06299212598425197
$ws.Range("S18").Value = 0.05511811023622047
$ws.Range("F19").Value = 0.01559633027522936
$ws.Range("H19").Value = 0.1908256880733945
$ws.Range("I19").Value = 0.07981651376146789
$ws.Range("J19").Value = 0.4073394495412844
$ws.Range("K19").Value = 0.1155963302752294
$ws.Range("M19").Value = 0.01651376146788991
$ws.Range("O19").Value = 0.0853211009174312
$ws.Range("S19").Value = 0.0889908256880734
